$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.132.15'
Set-TextValue 'E2' '  +0.59%  '
Set-TextValue 'D3' '1.655.29'
Set-TextValue 'E3' '  +0.21%  '
Set-TextValue 'E4' '  -0.21%  '
Set-TextValue 'D5' '217.85'
Set-TextValue 'E5' '  +0.52%  '
Set-TextValue 'D6' '0.5314'
Set-TextValue 'E6' '  +2.48%  '
Set-TextValue 'E7' '  -0.20%  '
Set-TextValue 'E8' '  -0.21%  '
Set-TextValue 'D9' '0.06330'
Set-TextValue 'E9' '  +1.46%  '
Set-TextValue 'D10' '20.44'
Set-TextValue 'E10' '  -0.38%  '
Set-TextValue 'D11' '0.07793'
Set-TextValue 'E11' '  +0.96%  '
Set-TextValue 'D12' '4.524'
Set-TextValue 'E12' '  +1.51%  '
Set-TextValue 'D13' '1.682.56'
Set-TextValue 'E13' '  +1.83%  '
Set-TextValue 'D14' '1.882.26'
Set-TextValue 'E14' '  +0.02%  '
Set-TextValue 'D15' '0.5490'
Set-TextValue 'E15' '  +1.44%  '
Set-TextValue 'D16' '0.0₅8214'
Set-TextValue 'E16' '  +1.87%  '
Set-TextValue 'D17' '65.44'
Set-TextValue 'E17' '  +1.21%  '
Set-TextValue 'D18' '26.119.67'
Set-TextValue 'E18' '  +0.39%  '
Set-TextValue 'E19' '  -0.07%  '
Set-TextValue 'D20' '4.603'
Set-TextValue 'E20' '  +0.91%  '
Set-TextValue 'D21' '191.13'
Set-TextValue 'E22' '  +0.99%  '
Set-TextValue 'D23' '6.026'
Set-TextValue 'E23' '  +1.00%  '
Set-TextValue 'E24' '  -0.12%  '
Set-TextValue 'D25' '145.42'
Set-TextValue 'E25' '  +5.43%  '
Set-TextValue 'D26' '0.1230'
Set-TextValue 'E26' '  +0.09%  '
Set-TextValue 'D27' '7.219'
Set-TextValue 'E27' '  +0.10%  '
Set-TextValue 'E28' '  -0.33%  '
Set-TextValue 'D29' '1.461'
Set-TextValue 'E29' '  +4.24%  '
Set-TextValue 'D30' '0.05790'
Set-TextValue 'E30' '  -1.86%  '
Set-TextValue 'E31' '  +0.00%  '
Set-TextValue 'D32' '3.559'
Set-TextValue 'E32' '  +1.28%  '
Set-TextValue 'D33' '3.276'
Set-TextValue 'E33' '  +1.01%  '
Set-TextValue 'D34' '1.604'
Set-TextValue 'E34' '  +2.70%  '
Set-TextValue 'E35' '  +1.53%  '
Set-TextValue 'D36' '0.9517'
Set-TextValue 'E36' '  +0.61%  '
Set-TextValue 'D37' '2.416'
Set-TextValue 'E37' '  -0.04%  '
Set-TextValue 'D38' '0.5759'
Set-TextValue 'E38' '  +2.32%  '
Set-TextValue 'D39' '0.01610'
Set-TextValue 'E39' '  +1.41%  '
Set-TextValue 'D40' '0.8560'
Set-TextValue 'E40' '  +1.18%  '
Set-TextValue 'D41' '5.793'
Set-TextValue 'E41' '  -1.53%  '
Set-TextValue 'D42' '104.74'
Set-TextValue 'E42' '  +3.96%  '
Set-TextValue 'E43' '  -0.13%  '
Set-TextValue 'D44' '1.038.38'
Set-TextValue 'E44' '  +4.02%  '
Set-TextValue 'D45' '1.796.47'
Set-TextValue 'E45' '  -0.08%  '
Set-TextValue 'D46' '57.00'
Set-TextValue 'E46' '  +1.21%  '
Set-TextValue 'D47' '1.005'
Set-TextValue 'E47' '  +0.63%  '
Set-TextValue 'D48' '0.4335'
Set-TextValue 'E48' '  +0.58%  '
Set-TextValue 'D49' '7.864'
Set-TextValue 'E49' '  -1.42%  '
Set-TextValue 'D50' '0.05149'
Set-TextValue 'E50' '  -0.02%  '
Set-TextValue 'D51' '1.447'
Set-TextValue 'E51' '  -1.06%  '
